$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of hour tracking data
$ws.Range("B27").Value = "Programmazione"
$ws.Range("C27").Value = "Creazione sistema salvataggio/load"
$ws.Range("B28").Value = "UI"
$ws.Range("C28").Value = "Creazione menu iniziale"

$ws.Range("D27").Value = 1/24
$ws.Range("E27").Value = "Per lo più preso da chatGPT, sob"

$ws.Range("D28").Value = 45/60/24
$ws.Range("E28").Value = "funge, manca però da capire come chiamare il load"

# Column C widens to fit the longer new text (bestFit)
$ws.Columns("C").ColumnWidth = 32

$ws.Range("E31").Select()
